$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.632.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.739.82"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4954"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06268"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.744.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07045"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6143"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.636.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007253"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.965.87"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.556"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.733"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.296"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.411"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.019"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08028"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.722"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04610"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9990"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.613"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.017"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6374"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9101"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.051"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.429"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01506"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.83"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.461"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3938"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.855"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1177"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05382"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.802"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.255"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.63%  "
